$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6,8).Value = 659.2727
$ws.Cells.Item(6,9).Value = 724.2
$ws.Cells.Item(6,10).Value = 10
$ws.Cells.Item(6,11).Value = 2172.6
$ws.Cells.Item(6,12).Value = 30
$ws.Cells.Item(6,13).Value = -2060.6
$ws.Cells.Item(6,14).Value = -254
$ws.Cells.Item(17,8).Value = 2494.2112
$ws.Cells.Item(17,9).Value = 0
$ws.Cells.Item(17,10).Value = 2494.2112
$ws.Cells.Item(17,11).Value = 0
$ws.Cells.Item(17,12).Value = 7482.633600000001
$ws.Cells.Item(17,13).ClearContents()
$ws.Cells.Item(17,14).Value = -7818.633600000001
$ws.Cells.Item(19,8).Value = 1457
$ws.Cells.Item(19,10).Value = 1666.2941
$ws.Cells.Item(19,12).Value = 1666.2941
$ws.Cells.Item(19,14).Value = -2016.2941
$ws.Cells.Item(32,8).Value = 6346.4614
$ws.Cells.Item(32,9).Value = 3133
$ws.Cells.Item(32,11).Value = 3133
$ws.Cells.Item(32,13).Value = -2807
$ws.Cells.Item(61,8).Value = 250
$ws.Cells.Item(61,9).Value = 250
$ws.Cells.Item(61,11).Value = 750
$ws.Cells.Item(61,13).Value = -578
$ws.Cells.Item(64,8).Value = 6372.143
$ws.Cells.Item(64,10).Value = 7446.923
$ws.Cells.Item(64,12).Value = 7446.923
$ws.Cells.Item(64,14).Value = -7942.923
$ws.Cells.Item(67,8).Value = 6372.143
$ws.Cells.Item(67,10).Value = 7446.923
$ws.Cells.Item(67,12).Value = 7446.923
$ws.Cells.Item(67,14).Value = -9162.922999999999
$ws.Cells.Item(68,8).Value = 0
$ws.Cells.Item(68,10).Value = 0
$ws.Cells.Item(68,12).Value = 0
$ws.Cells.Item(68,14).ClearContents()
$ws.Cells.Item(70,8).Value = 1226720.9
$ws.Cells.Item(70,9).Value = 1964
$ws.Cells.Item(70,11).Value = 5892
$ws.Cells.Item(70,13).Value = -5622
$ws.Cells.Item(71,8).Value = 0
$ws.Cells.Item(71,10).Value = 0
$ws.Cells.Item(71,12).Value = 0
$ws.Cells.Item(71,14).ClearContents()
$ws.Cells.Item(73,8).Value = 1226720.9
$ws.Cells.Item(73,9).Value = 1964
$ws.Cells.Item(73,11).Value = 5892
$ws.Cells.Item(73,13).Value = -4956
$ws.Cells.Item(86,8).Value = 4239.8
$ws.Cells.Item(86,9).Value = 2000
$ws.Cells.Item(86,11).Value = 2000
$ws.Cells.Item(86,13).Value = -877
$ws.Cells.Item(89,8).Value = 4239.8
$ws.Cells.Item(89,9).Value = 2000
$ws.Cells.Item(89,11).Value = 10000
$ws.Cells.Item(89,13).Value = -4384
$ws.Cells.Item(98,8).Value = 2487.2856
$ws.Cells.Item(98,9).Value = 2301.6667
$ws.Cells.Item(98,10).Value = 2626.5
$ws.Cells.Item(98,11).Value = 2301.6667
$ws.Cells.Item(98,12).Value = 2626.5
$ws.Cells.Item(98,13).Value = -803.6667000000002
$ws.Cells.Item(98,14).Value = -5622.5
$ws.Cells.Item(100,8).Value = 3382.4375
$ws.Cells.Item(100,9).Value = 978.3333
$ws.Cells.Item(100,10).Value = 4824.9
$ws.Cells.Item(100,11).Value = 978.3333
$ws.Cells.Item(100,12).Value = 4824.9
$ws.Cells.Item(100,13).Value = -437.3333
$ws.Cells.Item(100,14).Value = -5906.9
$ws.Cells.Item(106,8).Value = 2652.353
$ws.Cells.Item(106,9).Value = 2353.8462
$ws.Cells.Item(106,11).Value = 2353.8462
$ws.Cells.Item(106,13).Value = -1722.8462
$ws.Cells.Item(107,8).Value = 20834504
$ws.Cells.Item(107,10).Value = 12000
$ws.Cells.Item(107,12).Value = 12000
$ws.Cells.Item(107,14).Value = -15840
$ws.Cells.Item(111,8).Value = 566.38464
$ws.Cells.Item(111,9).Value = 605.9167
$ws.Cells.Item(111,11).Value = 1817.7501
$ws.Cells.Item(111,13).Value = 1249.2499
$ws.Cells.Item(113,8).Value = 252502.5
$ws.Cells.Item(113,9).Value = 3002
$ws.Cells.Item(113,10).Value = 502003
$ws.Cells.Item(113,11).Value = 3002
$ws.Cells.Item(113,12).Value = 502003
$ws.Cells.Item(113,13).Value = 252
$ws.Cells.Item(113,14).Value = -508511
$ws.Cells.Item(121,8).Value = 2583
$ws.Cells.Item(121,10).Value = 2583
$ws.Cells.Item(121,12).Value = 7749
$ws.Cells.Item(121,14).Value = -11243
$ws.Cells.Item(122,8).Value = 2487.2856
$ws.Cells.Item(122,9).Value = 2301.6667
$ws.Cells.Item(122,10).Value = 2626.5
$ws.Cells.Item(122,11).Value = 6905.000100000001
$ws.Cells.Item(122,12).Value = 7879.5
$ws.Cells.Item(122,13).Value = -4455.000100000001
$ws.Cells.Item(122,14).Value = -12779.5
$ws.Cells.Item(132,8).Value = 822.7368
$ws.Cells.Item(132,9).Value = 699.9423
$ws.Cells.Item(132,11).Value = 2099.8269
$ws.Cells.Item(132,13).Value = 430.1731
$ws.Cells.Item(135,8).Value = 974.7059
$ws.Cells.Item(135,9).Value = 970.875
$ws.Cells.Item(135,11).Value = 8737.875
$ws.Cells.Item(135,13).Value = -6202.875
$ws.Cells.Item(137,8).Value = 5779
$ws.Cells.Item(137,9).Value = 3510.84
$ws.Cells.Item(137,11).Value = 10532.52
$ws.Cells.Item(137,13).Value = -7982.52
$ws.Cells.Item(141,8).Value = 842.087
$ws.Cells.Item(141,9).Value = 842.087
$ws.Cells.Item(141,11).Value = 2526.261
$ws.Cells.Item(141,13).Value = 2653.739

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2,8).Value = 3540224.5
$ws.Cells.Item(2,9).Value = 4382273
$ws.Cells.Item(2,10).Value = 3622.2
$ws.Cells.Item(2,11).Value = 4382273
$ws.Cells.Item(2,12).Value = 3622.2
$ws.Cells.Item(2,13).Value = -4382160
$ws.Cells.Item(2,14).Value = -3848.2
$ws.Cells.Item(5,8).Value = 261.875
$ws.Cells.Item(5,9).Value = 261.875
$ws.Cells.Item(5,11).Value = 261.875
$ws.Cells.Item(5,13).Value = -149.875
$ws.Cells.Item(10,8).Value = 1333532.4
$ws.Cells.Item(10,9).Value = 2000149
$ws.Cells.Item(10,11).Value = 2000149
$ws.Cells.Item(10,13).Value = -1999979
$ws.Cells.Item(30,8).Value = 1130.1666
$ws.Cells.Item(30,9).Value = 1350.25
$ws.Cells.Item(30,10).Value = 690
$ws.Cells.Item(30,11).Value = 1350.25
$ws.Cells.Item(30,12).Value = 690
$ws.Cells.Item(30,13).Value = -1200.25
$ws.Cells.Item(30,14).Value = -990
$ws.Cells.Item(32,8).Value = 2945535.2
$ws.Cells.Item(32,9).Value = 3034294
$ws.Cells.Item(32,11).Value = 3034294
$ws.Cells.Item(32,13).Value = -3034007
$ws.Cells.Item(40,8).Value = 11166.333
$ws.Cells.Item(40,9).Value = 9249.5
$ws.Cells.Item(40,10).Value = 15000
$ws.Cells.Item(40,11).Value = 9249.5
$ws.Cells.Item(40,12).Value = 15000
$ws.Cells.Item(40,13).Value = -9073.5
$ws.Cells.Item(40,14).Value = -15352
$ws.Cells.Item(45,8).Value = 1737.6666
$ws.Cells.Item(45,9).Value = 1229.6666
$ws.Cells.Item(45,11).Value = 1229.6666
$ws.Cells.Item(45,13).Value = -852.6666
$ws.Cells.Item(61,8).Value = 3809.4822
$ws.Cells.Item(61,9).Value = 3140.423
$ws.Cells.Item(61,10).Value = 12507.25
$ws.Cells.Item(61,11).Value = 3140.423
$ws.Cells.Item(61,12).Value = 12507.25
$ws.Cells.Item(61,13).Value = -2928.423
$ws.Cells.Item(61,14).Value = -12931.25
$ws.Cells.Item(74,8).Value = 265578.97
$ws.Cells.Item(74,9).Value = 346222.7
$ws.Cells.Item(74,10).Value = 5727
$ws.Cells.Item(74,11).Value = 346222.7
$ws.Cells.Item(74,12).Value = 5727
$ws.Cells.Item(74,13).Value = -345348.7
$ws.Cells.Item(74,14).Value = -7475
$ws.Cells.Item(77,8).Value = 265578.97
$ws.Cells.Item(77,9).Value = 346222.7
$ws.Cells.Item(77,10).Value = 5727
$ws.Cells.Item(77,11).Value = 1731113.5
$ws.Cells.Item(77,12).Value = 28635
$ws.Cells.Item(77,13).Value = -1726745.5
$ws.Cells.Item(77,14).Value = -37371
$ws.Cells.Item(95,8).Value = 18600
$ws.Cells.Item(95,10).Value = 18750
$ws.Cells.Item(95,12).Value = 18750
$ws.Cells.Item(95,14).Value = -24242
$ws.Cells.Item(116,8).Value = 3540224.5
$ws.Cells.Item(116,9).Value = 4382273
$ws.Cells.Item(116,10).Value = 3622.2
$ws.Cells.Item(116,11).Value = 4382273
$ws.Cells.Item(116,12).Value = 3622.2
$ws.Cells.Item(116,13).Value = -4379979
$ws.Cells.Item(116,14).Value = -8210.200000000001
$ws.Cells.Item(122,8).Value = 2282.3333
$ws.Cells.Item(122,9).Value = 2490
$ws.Cells.Item(122,10).Value = 2178.5
$ws.Cells.Item(122,11).Value = 7470
$ws.Cells.Item(122,12).Value = 6535.5
$ws.Cells.Item(122,13).Value = -5020
$ws.Cells.Item(122,14).Value = -11435.5
$ws.Cells.Item(132,8).Value = 3902.1897
$ws.Cells.Item(132,9).Value = 3226.463
$ws.Cells.Item(132,10).Value = 13024.5
$ws.Cells.Item(132,11).Value = 9679.389000000001
$ws.Cells.Item(132,12).Value = 39073.5
$ws.Cells.Item(132,13).Value = -7149.389000000001
$ws.Cells.Item(132,14).Value = -44133.5
$ws.Cells.Item(136,8).Value = 3809.4822
$ws.Cells.Item(136,9).Value = 3140.423
$ws.Cells.Item(136,10).Value = 12507.25
$ws.Cells.Item(136,11).Value = 9421.269
$ws.Cells.Item(136,12).Value = 37521.75
$ws.Cells.Item(136,13).Value = -6871.269
$ws.Cells.Item(136,14).Value = -42621.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3,8).Value = 3540224.5
$ws.Cells.Item(3,9).Value = 4382273
$ws.Cells.Item(3,10).Value = 3622.2
$ws.Cells.Item(3,11).Value = 4382273
$ws.Cells.Item(3,12).Value = 3622.2
$ws.Cells.Item(3,13).Value = -4382159
$ws.Cells.Item(3,14).Value = -3850.2
$ws.Cells.Item(4,8).Value = 261.875
$ws.Cells.Item(4,9).Value = 261.875
$ws.Cells.Item(4,11).Value = 261.875
$ws.Cells.Item(4,13).Value = -146.875
$ws.Cells.Item(20,8).Value = 3755.6086
$ws.Cells.Item(20,9).Value = 3287.8
$ws.Cells.Item(20,10).Value = 4632.75
$ws.Cells.Item(20,11).Value = 3287.8
$ws.Cells.Item(20,12).Value = 4632.75
$ws.Cells.Item(20,13).Value = -3040.8
$ws.Cells.Item(20,14).Value = -5126.75
$ws.Cells.Item(134,8).Value = 6284.8184
$ws.Cells.Item(134,9).Value = 4361.636
$ws.Cells.Item(134,11).Value = 13084.908
$ws.Cells.Item(134,13).Value = -10549.908

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(11,8).Value = 1141.6
$ws.Cells.Item(11,9).Value = 782.6087
$ws.Cells.Item(11,10).Value = 2321.1428
$ws.Cells.Item(11,11).Value = 782.6087
$ws.Cells.Item(11,12).Value = 2321.1428
$ws.Cells.Item(11,13).Value = -642.6087
$ws.Cells.Item(11,14).Value = -2601.1428
$ws.Cells.Item(16,8).Value = 3946
$ws.Cells.Item(16,9).Value = 2931
$ws.Cells.Item(16,10).Value = 4816
$ws.Cells.Item(16,11).Value = 2931
$ws.Cells.Item(16,12).Value = 4816
$ws.Cells.Item(16,13).Value = -2644
$ws.Cells.Item(16,14).Value = -5390
$ws.Cells.Item(22,8).Value = 2472.1428
$ws.Cells.Item(22,9).Value = 372.85715
$ws.Cells.Item(22,11).Value = 372.85715
$ws.Cells.Item(22,13).Value = -22.85714999999999
$ws.Cells.Item(23,8).Value = 5087.5
$ws.Cells.Item(23,9).Value = 175
$ws.Cells.Item(23,11).Value = 175
$ws.Cells.Item(23,13).Value = 65
$ws.Cells.Item(27,8).Value = 5087.5
$ws.Cells.Item(27,9).Value = 175
$ws.Cells.Item(27,11).Value = 175
$ws.Cells.Item(27,13).Value = 17
$ws.Cells.Item(31,8).Value = 23813492
$ws.Cells.Item(31,9).Value = 62501572
$ws.Cells.Item(31,10).Value = 5445.115
$ws.Cells.Item(31,11).Value = 62501572
$ws.Cells.Item(31,12).Value = 5445.115
$ws.Cells.Item(31,13).Value = -62501277
$ws.Cells.Item(31,14).Value = -6035.115
$ws.Cells.Item(34,8).Value = 23813492
$ws.Cells.Item(34,9).Value = 62501572
$ws.Cells.Item(34,10).Value = 5445.115
$ws.Cells.Item(34,11).Value = 62501572
$ws.Cells.Item(34,12).Value = 5445.115
$ws.Cells.Item(34,13).Value = -62501370
$ws.Cells.Item(34,14).Value = -5849.115
$ws.Cells.Item(58,8).Value = 4381.766
$ws.Cells.Item(58,9).Value = 4281.3335
$ws.Cells.Item(58,10).Value = 4517.35
$ws.Cells.Item(58,11).Value = 4281.3335
$ws.Cells.Item(58,12).Value = 4517.35
$ws.Cells.Item(58,13).Value = -4078.3335
$ws.Cells.Item(58,14).Value = -4923.35
$ws.Cells.Item(62,8).Value = 15986.8125
$ws.Cells.Item(62,10).Value = 23430.5
$ws.Cells.Item(62,12).Value = 23430.5
$ws.Cells.Item(62,14).Value = -24678.5
$ws.Cells.Item(65,8).Value = 15986.8125
$ws.Cells.Item(65,10).Value = 23430.5
$ws.Cells.Item(65,12).Value = 117152.5
$ws.Cells.Item(65,14).Value = -123392.5
$ws.Cells.Item(86,8).Value = 9117.619000000001
$ws.Cells.Item(86,9).Value = 9279.4375
$ws.Cells.Item(86,11).Value = 9279.4375
$ws.Cells.Item(86,13).Value = -8156.4375
$ws.Cells.Item(89,8).Value = 9117.619000000001
$ws.Cells.Item(89,9).Value = 9279.4375
$ws.Cells.Item(89,11).Value = 46397.1875
$ws.Cells.Item(89,13).Value = -40781.1875
$ws.Cells.Item(105,8).Value = 912
$ws.Cells.Item(105,9).Value = 813.1
$ws.Cells.Item(105,11).Value = 813.1
$ws.Cells.Item(105,13).Value = 933.9
$ws.Cells.Item(107,8).Value = 659.8
$ws.Cells.Item(107,9).Value = 249
$ws.Cells.Item(107,11).Value = 249
$ws.Cells.Item(107,13).Value = 1671
$ws.Cells.Item(113,8).Value = 3946
$ws.Cells.Item(113,9).Value = 2931
$ws.Cells.Item(113,10).Value = 4816
$ws.Cells.Item(113,11).Value = 2931
$ws.Cells.Item(113,12).Value = 4816
$ws.Cells.Item(113,13).Value = -761
$ws.Cells.Item(113,14).Value = -9156
$ws.Cells.Item(132,8).Value = 5998.421
$ws.Cells.Item(132,9).Value = 5144.25
$ws.Cells.Item(132,10).Value = 7462.7144
$ws.Cells.Item(132,11).Value = 15432.75
$ws.Cells.Item(132,12).Value = 22388.1432
$ws.Cells.Item(132,13).Value = -12902.75
$ws.Cells.Item(132,14).Value = -27448.1432
$ws.Cells.Item(134,8).Value = 4263.485
$ws.Cells.Item(134,9).Value = 3934.1785
$ws.Cells.Item(134,11).Value = 11802.5355
$ws.Cells.Item(134,13).Value = -9267.5355
$ws.Cells.Item(136,8).Value = 4381.766
$ws.Cells.Item(136,9).Value = 4281.3335
$ws.Cells.Item(136,10).Value = 4517.35
$ws.Cells.Item(136,11).Value = 12844.0005
$ws.Cells.Item(136,12).Value = 13552.05
$ws.Cells.Item(136,13).Value = -10294.0005
$ws.Cells.Item(136,14).Value = -18652.05

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34,8).Value = 2686.6667
$ws.Cells.Item(34,9).Value = 194.28572
$ws.Cells.Item(34,10).Value = 4272.727
$ws.Cells.Item(34,11).Value = 582.85716
$ws.Cells.Item(34,12).Value = 12818.181
$ws.Cells.Item(34,13).Value = -498.85716
$ws.Cells.Item(34,14).Value = -12986.181
$ws.Cells.Item(36,8).Value = 489.4
$ws.Cells.Item(36,9).Value = 489.4
$ws.Cells.Item(36,11).Value = 1468.2
$ws.Cells.Item(36,13).Value = -1299.2
$ws.Cells.Item(54,8).Value = 2375
$ws.Cells.Item(54,10).Value = 4000
$ws.Cells.Item(54,12).Value = 12000
$ws.Cells.Item(54,14).Value = -13118
$ws.Cells.Item(64,8).Value = 2467.8
$ws.Cells.Item(64,9).Value = 1856.25
$ws.Cells.Item(64,10).Value = 4914
$ws.Cells.Item(64,11).Value = 5568.75
$ws.Cells.Item(64,12).Value = 14742
$ws.Cells.Item(64,13).Value = -5298.75
$ws.Cells.Item(64,14).Value = -15282
$ws.Cells.Item(67,8).Value = 2467.8
$ws.Cells.Item(67,9).Value = 1856.25
$ws.Cells.Item(67,10).Value = 4914
$ws.Cells.Item(67,11).Value = 5568.75
$ws.Cells.Item(67,12).Value = 14742
$ws.Cells.Item(67,13).Value = -4632.75
$ws.Cells.Item(67,14).Value = -16614
$ws.Cells.Item(87,8).Value = 0
$ws.Cells.Item(87,9).Value = 0
$ws.Cells.Item(87,11).Value = 0
$ws.Cells.Item(87,13).ClearContents()
$ws.Cells.Item(90,8).Value = 0
$ws.Cells.Item(90,9).Value = 0
$ws.Cells.Item(90,11).Value = 0
$ws.Cells.Item(90,13).ClearContents()
$ws.Cells.Item(103,8).Value = 2564.6667
$ws.Cells.Item(103,9).Value = 130
$ws.Cells.Item(103,10).Value = 4999.3335
$ws.Cells.Item(103,11).Value = 390
$ws.Cells.Item(103,12).Value = 14998.0005
$ws.Cells.Item(103,13).Value = 489
$ws.Cells.Item(103,14).Value = -16756.0005
$ws.Cells.Item(117,8).Value = 511.33334
$ws.Cells.Item(117,9).Value = 325.875
$ws.Cells.Item(117,11).Value = 977.625
$ws.Cells.Item(117,13).Value = 2464.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(29,8).Value = 333.33334
$ws.Cells.Item(29,9).Value = 325
$ws.Cells.Item(29,11).Value = 325
$ws.Cells.Item(29,13).Value = -35
$ws.Cells.Item(36,8).Value = 5704
$ws.Cells.Item(36,10).Value = 7500
$ws.Cells.Item(36,12).Value = 7500
$ws.Cells.Item(36,14).Value = -8470
$ws.Cells.Item(38,8).Value = 0
$ws.Cells.Item(38,10).Value = 0
$ws.Cells.Item(38,12).Value = 0
$ws.Cells.Item(38,14).ClearContents()
$ws.Cells.Item(70,8).Value = 6999.273
$ws.Cells.Item(70,9).Value = 6013.7144
$ws.Cells.Item(70,11).Value = 6013.7144
$ws.Cells.Item(70,13).Value = -5743.7144
$ws.Cells.Item(73,8).Value = 6999.273
$ws.Cells.Item(73,9).Value = 6013.7144
$ws.Cells.Item(73,11).Value = 6013.7144
$ws.Cells.Item(73,13).Value = -5077.7144
$ws.Cells.Item(97,8).Value = 1783.2954
$ws.Cells.Item(97,9).Value = 1588
$ws.Cells.Item(97,10).Value = 2815.5715
$ws.Cells.Item(97,11).Value = 1588
$ws.Cells.Item(97,12).Value = 2815.5715
$ws.Cells.Item(97,13).Value = -1092
$ws.Cells.Item(97,14).Value = -3807.5715
$ws.Cells.Item(102,8).Value = 1253.2821
$ws.Cells.Item(102,9).Value = 1137.5483
$ws.Cells.Item(102,11).Value = 1137.5483
$ws.Cells.Item(102,13).Value = 484.4517000000001
$ws.Cells.Item(107,8).Value = 1000
$ws.Cells.Item(107,9).Value = 1000
$ws.Cells.Item(107,10).Value = 0
$ws.Cells.Item(107,11).Value = 1000
$ws.Cells.Item(107,12).Value = 0
$ws.Cells.Item(107,13).Value = 920
$ws.Cells.Item(107,14).ClearContents()
$ws.Cells.Item(113,8).Value = 32292
$ws.Cells.Item(113,9).Value = 39302.5
$ws.Cells.Item(113,11).Value = 39302.5
$ws.Cells.Item(113,13).Value = -37132.5
$ws.Cells.Item(120,8).Value = 0
$ws.Cells.Item(120,10).Value = 0
$ws.Cells.Item(120,12).Value = 0
$ws.Cells.Item(120,14).ClearContents()
$ws.Cells.Item(122,8).Value = 2672.6304
$ws.Cells.Item(122,9).Value = 2689.7058
$ws.Cells.Item(122,11).Value = 8069.117400000001
$ws.Cells.Item(122,13).Value = -5619.117400000001
$ws.Cells.Item(126,8).Value = 4546
$ws.Cells.Item(126,9).Value = 4012
$ws.Cells.Item(126,11).Value = 12036
$ws.Cells.Item(126,13).Value = -9566
$ws.Cells.Item(132,8).Value = 5243.032
$ws.Cells.Item(132,9).Value = 3101.5908
$ws.Cells.Item(132,10).Value = 10477.667
$ws.Cells.Item(132,11).Value = 9304.7724
$ws.Cells.Item(132,12).Value = 31433.001
$ws.Cells.Item(132,13).Value = -6774.7724
$ws.Cells.Item(132,14).Value = -36493.001
$ws.Cells.Item(136,8).Value = 8131.727
$ws.Cells.Item(136,10).Value = 8131.727
$ws.Cells.Item(136,12).Value = 24395.181
$ws.Cells.Item(136,14).Value = -29495.181

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(35,8).Value = 1723.2
$ws.Cells.Item(35,9).Value = 1904.25
$ws.Cells.Item(35,10).Value = 999
$ws.Cells.Item(35,11).Value = 1904.25
$ws.Cells.Item(35,12).Value = 999
$ws.Cells.Item(35,13).Value = -1568.25
$ws.Cells.Item(35,14).Value = -1671
$ws.Cells.Item(43,8).Value = 13450
$ws.Cells.Item(43,10).Value = 21625
$ws.Cells.Item(43,12).Value = 21625
$ws.Cells.Item(43,14).Value = -22011
$ws.Cells.Item(55,8).Value = 772.75
$ws.Cells.Item(55,9).Value = 146.6
$ws.Cells.Item(55,10).Value = 1816.3334
$ws.Cells.Item(55,11).Value = 146.6
$ws.Cells.Item(55,12).Value = 1816.3334
$ws.Cells.Item(55,13).Value = 26.40000000000001
$ws.Cells.Item(55,14).Value = -2162.3334
$ws.Cells.Item(61,8).Value = 1906.4
$ws.Cells.Item(61,9).Value = 2181.75
$ws.Cells.Item(61,11).Value = 2181.75
$ws.Cells.Item(61,13).Value = -1979.75
$ws.Cells.Item(68,8).Value = 6138.6665
$ws.Cells.Item(68,9).Value = 3083.3333
$ws.Cells.Item(68,11).Value = 3083.3333
$ws.Cells.Item(68,13).Value = -2334.3333
$ws.Cells.Item(71,8).Value = 6138.6665
$ws.Cells.Item(71,9).Value = 3083.3333
$ws.Cells.Item(71,11).Value = 15416.6665
$ws.Cells.Item(71,13).Value = -11672.6665
$ws.Cells.Item(93,8).Value = 2725
$ws.Cells.Item(93,9).Value = 2725
$ws.Cells.Item(93,11).Value = 2725
$ws.Cells.Item(93,13).Value = -1477
$ws.Cells.Item(100,8).Value = 10002143
$ws.Cells.Item(100,9).Value = 83334000
$ws.Cells.Item(100,10).Value = 2344.0908
$ws.Cells.Item(100,11).Value = 83334000
$ws.Cells.Item(100,12).Value = 2344.0908
$ws.Cells.Item(100,13).Value = -83333459
$ws.Cells.Item(100,14).Value = -3426.0908
$ws.Cells.Item(113,8).Value = 1906.4
$ws.Cells.Item(113,9).Value = 2181.75
$ws.Cells.Item(113,11).Value = 2181.75
$ws.Cells.Item(113,13).Value = -11.75
$ws.Cells.Item(122,8).Value = 38464676
$ws.Cells.Item(122,9).Value = 50002716
$ws.Cells.Item(122,11).Value = 150008148
$ws.Cells.Item(122,13).Value = -150005698
$ws.Cells.Item(131,8).Value = 100969
$ws.Cells.Item(131,10).Value = 100969
$ws.Cells.Item(131,12).Value = 100969
$ws.Cells.Item(131,14).Value = -111049
$ws.Cells.Item(132,8).Value = 7294.143
$ws.Cells.Item(132,9).Value = 1986.0625
$ws.Cells.Item(132,10).Value = 24280
$ws.Cells.Item(132,11).Value = 5958.1875
$ws.Cells.Item(132,12).Value = 72840
$ws.Cells.Item(132,13).Value = -3428.1875
$ws.Cells.Item(132,14).Value = -77900
$ws.Cells.Item(133,8).Value = 124443.664
$ws.Cells.Item(133,10).Value = 124443.664
$ws.Cells.Item(133,12).Value = 124443.664
$ws.Cells.Item(133,14).Value = -129503.664
$ws.Cells.Item(135,8).Value = 318936.44
$ws.Cells.Item(135,10).Value = 318936.44
$ws.Cells.Item(135,12).Value = 318936.44
$ws.Cells.Item(135,14).Value = -329076.44
$ws.Cells.Item(136,8).Value = 4918.3687
$ws.Cells.Item(136,9).Value = 4253.5713
$ws.Cells.Item(136,10).Value = 6779.8
$ws.Cells.Item(136,11).Value = 12760.7139
$ws.Cells.Item(136,12).Value = 20339.4
$ws.Cells.Item(136,13).Value = -10210.7139
$ws.Cells.Item(136,14).Value = -25439.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32,8).Value = 7499.5
$ws.Cells.Item(32,9).Value = 7499.5
$ws.Cells.Item(32,11).Value = 7499.5
$ws.Cells.Item(32,13).Value = -7182.5
$ws.Cells.Item(122,8).Value = 4800
$ws.Cells.Item(122,9).Value = 4800
$ws.Cells.Item(122,10).Value = 0
$ws.Cells.Item(122,11).Value = 14400
$ws.Cells.Item(122,12).Value = 0
$ws.Cells.Item(122,13).Value = -11950
$ws.Cells.Item(122,14).ClearContents()
$ws.Cells.Item(132,8).Value = 5403.8623
$ws.Cells.Item(132,9).Value = 4957.0625
$ws.Cells.Item(132,11).Value = 14871.1875
$ws.Cells.Item(132,13).Value = -12341.1875
$ws.Cells.Item(133,8).Value = 100510.2
$ws.Cells.Item(133,10).Value = 100510.2
$ws.Cells.Item(133,12).Value = 100510.2
$ws.Cells.Item(133,14).Value = -110630.2
$ws.Cells.Item(136,8).Value = 2676.147
$ws.Cells.Item(136,9).Value = 973.6799999999999
$ws.Cells.Item(136,11).Value = 2921.04
$ws.Cells.Item(136,13).Value = -371.04
